$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 249.375  # H2
$ws.Cells.Item(2, 9).Value = 232.5  # I2
$ws.Cells.Item(2, 11).Value = 232.5  # K2
$ws.Cells.Item(2, 13).Value = -119.5  # M2

$ws.Cells.Item(28, 8).Value = 85199.75  # H28
$ws.Cells.Item(28, 9).Value = 85199.75  # I28
$ws.Cells.Item(28, 11).Value = 85199.75  # K28
$ws.Cells.Item(28, 13).Value = -84714.75  # M28

$ws.Cells.Item(43, 8).Value = 2052.6316  # H43
$ws.Cells.Item(43, 9).Value = 2000  # I43
$ws.Cells.Item(43, 10).Value = 2333.3333  # J43
$ws.Cells.Item(43, 11).Value = 2000  # K43
$ws.Cells.Item(43, 12).Value = 2333.3333  # L43
$ws.Cells.Item(43, 13).Value = -1931  # M43
$ws.Cells.Item(43, 14).Value = -2471.3333  # N43

$ws.Cells.Item(76, 8).Value = 90915630  # H76
$ws.Cells.Item(76, 9).Value = 5997.5  # I76
$ws.Cells.Item(76, 10).Value = 142864000  # J76
$ws.Cells.Item(76, 11).Value = 5997.5  # K76
$ws.Cells.Item(76, 12).Value = 142864000  # L76
$ws.Cells.Item(76, 13).Value = -5682.5  # M76
$ws.Cells.Item(76, 14).Value = -142864630  # N76

$ws.Cells.Item(79, 8).Value = 90915630  # H79
$ws.Cells.Item(79, 9).Value = 5997.5  # I79
$ws.Cells.Item(79, 10).Value = 142864000  # J79
$ws.Cells.Item(79, 11).Value = 5997.5  # K79
$ws.Cells.Item(79, 12).Value = 142864000  # L79
$ws.Cells.Item(79, 13).Value = -4905.5  # M79
$ws.Cells.Item(79, 14).Value = -142866184  # N79

$ws.Cells.Item(112, 8).Value = 1747.0555  # H112
$ws.Cells.Item(112, 9).Value = 0  # I112
$ws.Cells.Item(112, 10).Value = 1747.0555  # J112
$ws.Cells.Item(112, 11).Value = 0  # K112
$ws.Cells.Item(112, 12).Value = 5241.166499999999  # L112
$ws.Cells.Item(112, 13).ClearContents()  # M112
$ws.Cells.Item(112, 14).Value = -7457.166499999999  # N112

$ws.Cells.Item(138, 8).Value = 5386.2446  # H138
$ws.Cells.Item(138, 9).Value = 2456.5144  # I138
$ws.Cells.Item(138, 10).Value = 7124.22  # J138
$ws.Cells.Item(138, 11).Value = 7369.5432  # K138
$ws.Cells.Item(138, 12).Value = 21372.66  # L138
$ws.Cells.Item(138, 13).Value = -2229.5432  # M138
$ws.Cells.Item(138, 14).Value = -31652.66  # N138

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 2422.1738  # H74
$ws.Cells.Item(74, 9).Value = 2511.0527  # I74
$ws.Cells.Item(74, 10).Value = 2000  # J74
$ws.Cells.Item(74, 11).Value = 2511.0527  # K74
$ws.Cells.Item(74, 12).Value = 2000  # L74
$ws.Cells.Item(74, 13).Value = -1637.0527  # M74
$ws.Cells.Item(74, 14).Value = -3748  # N74

$ws.Cells.Item(77, 8).Value = 2422.1738  # H77
$ws.Cells.Item(77, 9).Value = 2511.0527  # I77
$ws.Cells.Item(77, 10).Value = 2000  # J77
$ws.Cells.Item(77, 11).Value = 12555.2635  # K77
$ws.Cells.Item(77, 12).Value = 10000  # L77
$ws.Cells.Item(77, 13).Value = -8187.263500000001  # M77
$ws.Cells.Item(77, 14).Value = -18736  # N77

$ws.Cells.Item(88, 8).Value = 4697.0835  # H88
$ws.Cells.Item(88, 10).Value = 3820.875  # J88
$ws.Cells.Item(88, 12).Value = 3820.875  # L88
$ws.Cells.Item(88, 14).Value = -4632.875  # N88

$ws.Cells.Item(91, 8).Value = 4697.0835  # H91
$ws.Cells.Item(91, 10).Value = 3820.875  # J91
$ws.Cells.Item(91, 12).Value = 3820.875  # L91
$ws.Cells.Item(91, 14).Value = -6628.875  # N91

$ws.Cells.Item(109, 8).Value = 99950  # H109
$ws.Cells.Item(109, 10).Value = 99950  # J109
$ws.Cells.Item(109, 12).Value = 99950  # L109
$ws.Cells.Item(109, 14).Value = -102724  # N109

$ws.Cells.Item(112, 8).Value = 37000.332  # H112
$ws.Cells.Item(112, 10).Value = 38000.5  # J112
$ws.Cells.Item(112, 12).Value = 38000.5  # L112
$ws.Cells.Item(112, 14).Value = -40954.5  # N112

$ws.Cells.Item(122, 8).Value = 3664  # H122
$ws.Cells.Item(122, 9).Value = 2382.0908  # I122
$ws.Cells.Item(122, 10).Value = 5678.4287  # J122
$ws.Cells.Item(122, 11).Value = 7146.2724  # K122
$ws.Cells.Item(122, 12).Value = 17035.2861  # L122
$ws.Cells.Item(122, 13).Value = -4696.2724  # M122
$ws.Cells.Item(122, 14).Value = -21935.2861  # N122

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 207.14285  # H22
$ws.Cells.Item(22, 9).Value = 225  # I22
$ws.Cells.Item(22, 10).Value = 100  # J22
$ws.Cells.Item(22, 11).Value = 225  # K22
$ws.Cells.Item(22, 12).Value = 100  # L22
$ws.Cells.Item(22, 13).Value = -52  # M22
$ws.Cells.Item(22, 14).Value = -446  # N22

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2884.0557  # H58
$ws.Cells.Item(58, 9).Value = 2806.5  # I58
$ws.Cells.Item(58, 11).Value = 2806.5  # K58
$ws.Cells.Item(58, 13).Value = -2603.5  # M58

$ws.Cells.Item(99, 8).Value = 4473.7617  # H99
$ws.Cells.Item(99, 9).Value = 3643.2354  # I99
$ws.Cells.Item(99, 11).Value = 3643.2354  # K99
$ws.Cells.Item(99, 13).Value = -2145.2354  # M99

$ws.Cells.Item(107, 8).Value = 883.25  # H107
$ws.Cells.Item(107, 9).Value = 870.25  # I107
$ws.Cells.Item(107, 10).Value = 922.25  # J107
$ws.Cells.Item(107, 11).Value = 870.25  # K107
$ws.Cells.Item(107, 12).Value = 922.25  # L107
$ws.Cells.Item(107, 13).Value = 1049.75  # M107
$ws.Cells.Item(107, 14).Value = -4762.25  # N107

$ws.Cells.Item(126, 8).Value = 4473.7617  # H126
$ws.Cells.Item(126, 9).Value = 3643.2354  # I126
$ws.Cells.Item(126, 11).Value = 10929.7062  # K126
$ws.Cells.Item(126, 13).Value = -8459.706200000001  # M126

$ws.Cells.Item(132, 8).Value = 2226  # H132
$ws.Cells.Item(132, 9).Value = 1321.8  # I132
$ws.Cells.Item(132, 11).Value = 3965.4  # K132
$ws.Cells.Item(132, 13).Value = -1435.4  # M132

$ws.Cells.Item(134, 8).Value = 457474.53  # H134
$ws.Cells.Item(134, 9).Value = 3068.5715  # I134
$ws.Cells.Item(134, 10).Value = 10000000  # J134
$ws.Cells.Item(134, 11).Value = 9205.7145  # K134
$ws.Cells.Item(134, 12).Value = 30000000  # L134
$ws.Cells.Item(134, 13).Value = -6670.7145  # M134
$ws.Cells.Item(134, 14).Value = -30005070  # N134

$ws.Cells.Item(136, 8).Value = 2884.0557  # H136
$ws.Cells.Item(136, 9).Value = 2806.5  # I136
$ws.Cells.Item(136, 11).Value = 8419.5  # K136
$ws.Cells.Item(136, 13).Value = -5869.5  # M136

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(50, 8).Value = 38739.383  # H50
$ws.Cells.Item(50, 10).Value = 62629.125  # J50
$ws.Cells.Item(50, 12).Value = 187887.375  # L50
$ws.Cells.Item(50, 14).Value = -188849.375  # N50

$ws.Cells.Item(53, 8).Value = 38739.383  # H53
$ws.Cells.Item(53, 10).Value = 62629.125  # J53
$ws.Cells.Item(53, 12).Value = 187887.375  # L53
$ws.Cells.Item(53, 14).Value = -188849.375  # N53

$ws.Cells.Item(131, 8).Value = 2826.2407  # H131
$ws.Cells.Item(131, 9).Value = 2385.8  # I131
$ws.Cells.Item(131, 10).Value = 2871.1836  # J131
$ws.Cells.Item(131, 11).Value = 7157.400000000001  # K131
$ws.Cells.Item(131, 12).Value = 8613.550799999999  # L131
$ws.Cells.Item(131, 13).Value = -2117.400000000001  # M131
$ws.Cells.Item(131, 14).Value = -18693.5508  # N131

$ws.Cells.Item(134, 8).Value = 4573.25  # H134
$ws.Cells.Item(134, 9).Value = 4573.25  # I134
$ws.Cells.Item(134, 11).Value = 13719.75  # K134
$ws.Cells.Item(134, 13).Value = -8649.75  # M134

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(103, 8).Value = 48111  # H103
$ws.Cells.Item(103, 10).Value = 48111  # J103
$ws.Cells.Item(103, 12).Value = 48111  # L103
$ws.Cells.Item(103, 14).Value = -50455  # N103

$ws.Cells.Item(122, 8).Value = 3937.8948  # H122
$ws.Cells.Item(122, 9).Value = 3001.4285  # I122
$ws.Cells.Item(122, 11).Value = 9004.2855  # K122
$ws.Cells.Item(122, 13).Value = -6554.2855  # M122

$ws.Cells.Item(132, 8).Value = 41608.11  # H132
$ws.Cells.Item(132, 9).Value = 4976.0415  # I132
$ws.Cells.Item(132, 11).Value = 14928.1245  # K132
$ws.Cells.Item(132, 13).Value = -12398.1245  # M132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 7184.385  # H7
$ws.Cells.Item(7, 9).Value = 7300.778  # I7
$ws.Cells.Item(7, 11).Value = 7300.778  # K7
$ws.Cells.Item(7, 13).Value = -7188.778  # M7

$ws.Cells.Item(40, 8).Value = 3405030  # H40
$ws.Cells.Item(40, 9).Value = 7287350  # I40
$ws.Cells.Item(40, 11).Value = 7287350  # K40
$ws.Cells.Item(40, 13).Value = -7287214  # M40

$ws.Cells.Item(61, 8).Value = 2361.0667  # H61
$ws.Cells.Item(61, 9).Value = 2493.8076  # I61
$ws.Cells.Item(61, 11).Value = 2493.8076  # K61
$ws.Cells.Item(61, 13).Value = -2291.8076  # M61

$ws.Cells.Item(93, 8).Value = 2722.3257  # H93
$ws.Cells.Item(93, 9).Value = 2545.7334  # I93
$ws.Cells.Item(93, 10).Value = 3129.8462  # J93
$ws.Cells.Item(93, 11).Value = 2545.7334  # K93
$ws.Cells.Item(93, 12).Value = 3129.8462  # L93
$ws.Cells.Item(93, 13).Value = -1297.7334  # M93
$ws.Cells.Item(93, 14).Value = -5625.8462  # N93

$ws.Cells.Item(113, 8).Value = 2361.0667  # H113
$ws.Cells.Item(113, 9).Value = 2493.8076  # I113
$ws.Cells.Item(113, 11).Value = 2493.8076  # K113
$ws.Cells.Item(113, 13).Value = -323.8076000000001  # M113

$ws.Cells.Item(122, 8).Value = 1006599.2  # H122
$ws.Cells.Item(122, 9).Value = 3336664  # I122
$ws.Cells.Item(122, 11).Value = 10009992  # K122
$ws.Cells.Item(122, 13).Value = -10007542  # M122

$ws.Cells.Item(126, 8).Value = 7184.385  # H126
$ws.Cells.Item(126, 9).Value = 7300.778  # I126
$ws.Cells.Item(126, 11).Value = 21902.334  # K126
$ws.Cells.Item(126, 13).Value = -19432.334  # M126

$ws.Cells.Item(132, 8).Value = 3929  # H132
$ws.Cells.Item(132, 9).Value = 3442.2  # I132
$ws.Cells.Item(132, 11).Value = 10326.6  # K132
$ws.Cells.Item(132, 13).Value = -7796.599999999999  # M132

$ws.Cells.Item(136, 8).Value = 230682.58  # H136
$ws.Cells.Item(136, 9).Value = 351231.12  # I136
$ws.Cells.Item(136, 10).Value = 12188.3125  # J136
$ws.Cells.Item(136, 11).Value = 1053693.36  # K136
$ws.Cells.Item(136, 12).Value = 36564.9375  # L136
$ws.Cells.Item(136, 13).Value = -1051143.36  # M136
$ws.Cells.Item(136, 14).Value = -41664.9375  # N136

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(16, 8).Value = 114750  # H16
$ws.Cells.Item(16, 10).Value = 114750  # J16
$ws.Cells.Item(16, 12).Value = 114750  # L16
$ws.Cells.Item(16, 14).Value = -115334  # N16

$ws.Cells.Item(136, 8).Value = 368157.12  # H136
$ws.Cells.Item(136, 9).Value = 439035.53  # I136
$ws.Cells.Item(136, 10).Value = 205136.8  # J136
$ws.Cells.Item(136, 11).Value = 1317106.59  # K136
$ws.Cells.Item(136, 12).Value = 615410.3999999999  # L136
$ws.Cells.Item(136, 13).Value = -1314556.59  # M136
$ws.Cells.Item(136, 14).Value = -620510.3999999999  # N136
